$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New date row: append "27-10-2025" and the gold price text to row 43,
# mirroring the style/format used by the previous rows (A column = date, B column = price text).
$ws.Range("A43").Value = "27-10-2025"
$ws.Range("B43").Value = "The price of gold in India today is ₹12,448 per gram for 24 karat gold, ₹11,410 per gram for 22 karat gold and ₹9,336 per gram for 18 karat gold (also called 999 gold)."

$ws.Range("A42").Copy() | Out-Null
$ws.Range("A43").PasteSpecial(-4122) | Out-Null

$ws.Range("B42").Copy() | Out-Null
$ws.Range("B43").PasteSpecial(-4122) | Out-Null
